$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: clear the old (placeholder) entry ---
$ws.Range("A2").ClearContents()
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()

# --- Row 3: "ecouter le prof" / "Le prof a parlé sur Les MCD" ---
$ws.Range("A3").Value = 44245
$ws.Range("B3").Value = 0.33333333333333331
$ws.Range("C3").Value = 0.35416666666666669
$ws.Range("E3").Value = "Théorie"
$ws.Range("F3").Value = "ecouter le prof"
$ws.Range("G3").Value = "Travail"
$ws.Range("H3").Value = "MA-20"
$ws.Range("I3").Value = "Le prof a parlé sur Les MCD"
$ws.Range("J3").Value = "OUI"

# --- Row 4: "faire le MCD,Scénario," / "j'ai commencé le MCD" ---
$ws.Range("A4").Value = 44245
$ws.Range("B4").Value = 0.35416666666666669
$ws.Range("C4").Value = 0.36458333333333331
$ws.Range("E4").Value = "Pratique"
$ws.Range("F4").Value = "faire le MCD,Scénario,"
$ws.Range("G4").Value = "Travail"
$ws.Range("H4").Value = "MA-20"
$ws.Range("I4").Value = "j'ai commencé le MCD"
$ws.Range("J4").Value = "NON"

# --- Row 5: "j'ai continuer le MCD" ---
$ws.Range("A5").Value = 44245
$ws.Range("B5").Value = 0.36805555555555558
$ws.Range("C5").Value = 0.39930555555555558
$ws.Range("E5").Value = "Pratique"
$ws.Range("F5").Value = "faire le MCD,Scénario,"
$ws.Range("G5").Value = "Travail"
$ws.Range("H5").Value = "MA-20"
$ws.Range("I5").Value = "j'ai continuer le MCD"
$ws.Range("J5").Value = "NON"

# --- Row 6: "j'ai continuer le MCD et commencer le Scénario" ---
$ws.Range("A6").Value = 44245
$ws.Range("B6").Value = 0.40972222222222227
$ws.Range("C6").Value = 0.44097222222222227
$ws.Range("E6").Value = "Pratique"
$ws.Range("F6").Value = "faire le MCD,Scénario,"
$ws.Range("G6").Value = "Travail"
$ws.Range("H6").Value = "MA-20"
$ws.Range("I6").Value = "j'ai continuer le MCD et commencer le Scénario"
$ws.Range("J6").Value = "NON"

# --- Row 7: "j'ai continuer le MCD et continuer le Scénario" ---
$ws.Range("A7").Value = 44245
$ws.Range("B7").Value = 0.44444444444444442
$ws.Range("C7").Value = 0.47569444444444442
$ws.Range("E7").Value = "Pratique"
$ws.Range("F7").Value = "faire le MCD,Scénario,"
$ws.Range("G7").Value = "Travail"
$ws.Range("H7").Value = "MA-20"
$ws.Range("I7").Value = "j'ai continuer le MCD et continuer le Scénario"
$ws.Range("J7").Value = "NON"

# --- Row 8: same descriptif as row 7 ---
$ws.Range("A8").Value = 44245
$ws.Range("B8").Value = 0.47916666666666669
$ws.Range("C8").Value = 0.51041666666666663
$ws.Range("E8").Value = "Pratique"
$ws.Range("F8").Value = "faire le MCD,Scénario,"
$ws.Range("G8").Value = "Travail"
$ws.Range("H8").Value = "MA-20"
$ws.Range("I8").Value = "j'ai continuer le MCD et continuer le Scénario"
$ws.Range("J8").Value = "NON"

# --- Row 11 is authored before row 10 in the original session (shared-string
#     insertion order), so set it first to keep the sharedStrings table
#     byte-identical to the target ---
$ws.Range("A11").Value = 44246
$ws.Range("B11").Value = 0.59722222222222221
$ws.Range("C11").Value = 0.62847222222222221
$ws.Range("E11").Value = "Pratique"
$ws.Range("F11").Value = "faire le MCD,Scénario,"
$ws.Range("G11").Value = "Travail"
$ws.Range("H11").Value = "MA-20"
$ws.Range("I11").Value = "j'ai fini le MCD et Le scénario"
$ws.Range("J11").Value = "OUI"

# --- Row 10: "j'ai commencé le code de la bataille navale ,j'ai fait le menu" ---
$ws.Range("A10").Value = 44246
$ws.Range("B10").Value = 0.5625
$ws.Range("C10").Value = 0.59375
$ws.Range("E10").Value = "Pratique"
$ws.Range("F10").Value = "faire le MCD,Scénario,"
$ws.Range("G10").Value = "Travail"
$ws.Range("H10").Value = "MA-20"
$ws.Range("I10").Value = "j'ai commencé le code de la bataille navale ,j'ai fait le menu"
$ws.Range("J10").Value = "NON"

# --- Column I got a bit wider to fit the new, longer descriptions ---
$ws.Columns.Item(9).ColumnWidth = 51.498697916666664

# --- Move the active selection like the author left it ---
$ws.Range("A12").Select()
